# Update cryptocurrency price/volume data per commit
# "Updated cryptos list on Sat Apr  8 10:30:19 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "28.099.85"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.873.31"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'313.02"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5121"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.3904"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("D9").Value = "'0.08315"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "'41.48"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'6.203"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'20.63"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.876.30"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D18").Value = "'90.94"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "'0.06650"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "'17.72"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D23").Value = "28.139.85"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "'2.257"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'3.393"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.092.07"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.486"
$ws.Range("E28").Value = "  -3.28%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'159.07"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'20.59"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'124.86"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1062"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.037"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.848"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.606"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.608"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02449"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06557"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2182"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.199"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6477"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.227"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "'4.978"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'11.27"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6130"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.05"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.285"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.664"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.016"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.231"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'120.41"
$ws.Range("E51").Value = "  -0.22%  "
